$wb = $excel.ActiveWorkbook

# --- students sheet: enrollment counts were wrong (2 -> 14) ---
$ws1 = $wb.Worksheets.Item("students")
$ws1.Range("A2:A24").Value = 14

# --- semester sheet: fix year and remove the duplicated/incorrect row that
#     was losing enrollment table records ---
$ws2 = $wb.Worksheets.Item("semester")
$ws2.Range("A2").Value = 2015
$ws2.Rows.Item(4).Delete()

# --- Update which sheet/cell is active & selected ---
$ws1.Activate()
$ws1.Range("A25:XFD25").Select()

$ws2.Activate()
$ws2.Range("K16").Select()

$ws1.Activate()
